$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the right-hand panel area with a purple fill
$ws.Range("I2:M8").Interior.Color = 8388736

# "Console panel" label, merged across I4:M4, centered
$ws.Range("I4:M4").Merge()
$ws.Range("I4").Value = "Console panel"
$ws.Range("I4:M4").HorizontalAlignment = -4108

# "Gameplay panel" label, merged across C5:G5
$ws.Range("C5:G5").Merge()
$ws.Range("C5").Value = "Gameplay panel"

# Update view: zoom and active cell selection
$excel.ActiveWindow.Zoom = 140
$ws.Range("C6").Select()
